$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats = -4122, xlPasteValues = -4163 (Excel COM constants)

# --- Row 10 (Objetivos:) -------------------------------------------------
# The long "Complementar a formacao..." blurb is replaced by the
# "Docentes responsaveis" value (reuse text straight from B13/C13 so the
# shared-string + style stay byte-identical).
$ws.Range("B13").Copy()
$ws.Range("B10").PasteSpecial(-4163)
$ws.Range("C13").Copy()
$ws.Range("C10").PasteSpecial(-4163)

# --- Row 13 ---------------------------------------------------------------
# Previously held only B/C ("11079086 - ..."); now also gets the
# "Programa resumido:" label in column A, and the B/C text becomes
# "Semestral". Column A already defaults to style 1 for new cells here.
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

# --- Row 14: becomes "Short syllabus:" / "To be defined..." -------------
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "To be defined according to the scheduled topic"
$ws.Range("C14").Value = "To be defined according to the scheduled topic"
$ws.Rows(14).RowHeight = 60

# --- Row 15: becomes "Programa:" / "01/01/2021" ---------------------------
# Reuse the literal text value from B8/C8 (which already store
# "01/01/2021" as a shared string) via PasteSpecial-values so the engine
# doesn't reinterpret the literal as a date and spawn a new number format.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$ws.Rows(15).RowHeight = 120

# --- Row 16: becomes "Syllabus:" / long English syllabus text ------------
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "The content of this optional course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."
$ws.Range("C16").Value = "The content of this optional course will be according to the topic to be programmed, and should address complementary subjects to the regular content of the undergraduate course."
$ws.Rows(16).RowHeight = 120

# --- Row 17: becomes "Avaliacao:" label only; drop its B/C values --------
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Rows(17).EntireRow.AutoFit()

# --- Row 18: becomes "Metodo:" / "11079086 - Herlandí de Souza Andrade" --
# B18/C18 did not exist before, so first pull in the B10/C10 formatting
# (style 2 / style 3) before writing the value, otherwise the engine
# defaults new cells in this row to the bold label style.
$ws.Range("A18").Value = "Método:"
$ws.Range("B10").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C10").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Rows(18).RowHeight = 60

# --- Row 19: becomes "Criterio:" / evaluation paragraph ------------------
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."
$ws.Range("C19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."
$ws.Rows(19).RowHeight = 60

# --- Row 20: becomes "Norma de recuperacao:" / "Media ponderada..." ------
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média ponderada das avaliações (M)."
$ws.Range("C20").Value = "Média ponderada das avaliações (M)."
$ws.Rows(20).RowHeight = 60

# --- Row 21: becomes "Bibliografia:" / recovery paragraph ----------------
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
$ws.Range("C21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
$ws.Rows(21).RowHeight = 120

# --- Row 22: becomes "Requisitos:" label only; drop its B/C values -------
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Rows(22).EntireRow.AutoFit()

# --- Row 23: drop the "Requisitos:" label from A; move the requirement ---
# text (previously row 24) up into B/C here. Paste formats first (B23/C23
# did not exist before) then the literal values from B24/C24.
$ws.Range("A23").Clear()
$ws.Range("B24").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B24").Copy()
$ws.Range("B23").PasteSpecial(-4163)
$ws.Range("C24").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C24").Copy()
$ws.Range("C23").PasteSpecial(-4163)
$ws.Rows(23).RowHeight = 30

# --- Row 24 no longer exists; remove it entirely --------------------------
$ws.Rows(24).Delete()
